# Applies the commit's text edits to poster/Presentation1.pptx (slide 1):
#   1. "...並依序利用DP找出能轉至基底的路徑"
#        -> "...並依序利用DFS找出能轉至基底的路徑"   (DP -> DFS, "找出" split into its own run)
#   2. "...最後根據此陣列利用DP找出最佳縫合影像的順序"
#        -> "...最後根據此陣列利用DFS找出最佳縫合影像的順序" (DP -> DFS, "利用"/"找出" split off)
#   3. "根據縫合順序縫合影像" (was two runs "根據縫合順序縫合" + "影像")
#        -> merged into a single run with the same text
#
# Helper: replace the first occurrence of $needle (searched fresh each time,
# since prior edits can shift offsets) with $replacement, inside the given
# TextRange2. Using Characters(Start, Length) causes PowerPoint to split
# the underlying run(s) at the edit boundaries while leaving untouched runs
# alone, exactly like typing over a selection.
function Replace-FirstText {
    param(
        $TextRange,
        [string]$Needle,
        [string]$Replacement
    )
    $full = $TextRange.Text
    $idx0 = $full.IndexOf($Needle)
    if ($idx0 -lt 0) {
        throw "Replace-FirstText: '$Needle' not found"
    }
    $chars = $TextRange.Characters($idx0 + 1, $Needle.Length)
    $chars.Text = $Replacement
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---- Change 1 -----------------------------------------------------------
# Shape id=8 "TextBox 8", nested inside the group chain Shapes.Item(5).
$g1 = $s.Shapes.Item(5)
$tb1 = $g1.GroupItems.Item(1)
$tr1 = $tb1.TextFrame2.TextRange

# Split "DP" away from the following text so it becomes its own run, then
# retype it as "DFS".
Replace-FirstText $tr1 "DP" "DFS"
# Split "找出" off the front of "找出能轉至基底的路徑" into its own run
# (re-typed with identical text so only the run boundary changes).
Replace-FirstText $tr1 "找出" "找出"

# ---- Change 2 -------------------------------------------------------------
# Shape id=17 "TextBox 8", nested inside the group chain Shapes.Item(6).
$g2 = $s.Shapes.Item(6)
$tb2 = $g2.GroupItems.Item(1)
$tr2 = $tb2.TextFrame2.TextRange

# Split "利用" off the end of the long sentence run so it becomes its own run.
Replace-FirstText $tr2 "利用DP" "利用DP"
# DP -> DFS
Replace-FirstText $tr2 "DP" "DFS"
# Split "找出" off the front of "找出最佳縫合影像的順序".
Replace-FirstText $tr2 "找出" "找出"

# ---- Change 3 -------------------------------------------------------------
# Merge the two runs "根據縫合順序縫合" + "影像" into a single run by
# retyping the full combined span in one go.
Replace-FirstText $tr2 "根據縫合順序縫合影像" "根據縫合順序縫合影像"
